$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'68.680.59"
$ws.Cells.Item(2, 5).Value = '  -0.40%  '

$ws.Cells.Item(3, 4).Value = "'3.758.17"
$ws.Cells.Item(3, 5).Value = '  -1.49%  '

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

$ws.Cells.Item(5, 4).Value = "'628.24"
$ws.Cells.Item(5, 5).Value = '  +0.24%  '

$ws.Cells.Item(6, 4).Value = "'164.70"
$ws.Cells.Item(6, 5).Value = '  -0.17%  '

$ws.Cells.Item(7, 4).Value = "'3.755.83"
$ws.Cells.Item(7, 5).Value = '  -1.47%  '

$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$ws.Cells.Item(9, 5).Value = '  -0.05%  '

$ws.Cells.Item(10, 4).Value = "'0.157"
$ws.Cells.Item(10, 5).Value = '  -2.25%  '

$ws.Cells.Item(11, 5).Value = '  -0.06%  '

$ws.Cells.Item(12, 4).Value = "'6.90"
$ws.Cells.Item(12, 5).Value = '  +4.41%  '

$ws.Cells.Item(13, 4).Value = "'0.0000237"
$ws.Cells.Item(13, 5).Value = '  -4.95%  '

$ws.Cells.Item(14, 4).Value = "'34.72"
$ws.Cells.Item(14, 5).Value = '  -3.42%  '

$ws.Cells.Item(15, 4).Value = "'4.386.31"
$ws.Cells.Item(15, 5).Value = '  -1.55%  '

$ws.Cells.Item(16, 4).Value = "'3.749.40"
$ws.Cells.Item(16, 5).Value = '  -2.70%  '

$ws.Cells.Item(17, 4).Value = "'68.651.96"

$ws.Cells.Item(18, 4).Value = "'17.59"
$ws.Cells.Item(18, 5).Value = '  -2.39%  '

$ws.Cells.Item(19, 5).Value = '  -0.31%  '

$ws.Cells.Item(20, 4).Value = "'6.97"
$ws.Cells.Item(20, 5).Value = '  -1.99%  '

$ws.Cells.Item(21, 4).Value = "'465.89"
$ws.Cells.Item(21, 5).Value = '  +0.25%  '

$ws.Cells.Item(22, 4).Value = "'9.46"
$ws.Cells.Item(22, 5).Value = '  -2.03%  '

$ws.Cells.Item(23, 4).Value = "'0.700"
$ws.Cells.Item(23, 5).Value = '  -0.99%  '

$ws.Cells.Item(24, 4).Value = "'81.46"
$ws.Cells.Item(24, 5).Value = '  -2.65%  '

$ws.Cells.Item(25, 5).Value = '  -6.37%  '

$ws.Cells.Item(26, 4).Value = "'12.05"
$ws.Cells.Item(26, 5).Value = '  +1.23%  '

$ws.Cells.Item(27, 4).Value = "'2.09"
$ws.Cells.Item(27, 5).Value = '  -2.23%  '

$ws.Cells.Item(28, 4).Value = "'10.04"
$ws.Cells.Item(28, 5).Value = '  +0.35%  '

$ws.Cells.Item(29, 4).Value = "'1.00"
$ws.Cells.Item(29, 5).Value = '  -0.03%  '

$ws.Cells.Item(30, 4).Value = "'3.905.74"
$ws.Cells.Item(30, 5).Value = '  -1.54%  '

$ws.Cells.Item(31, 4).Value = "'2.26"
$ws.Cells.Item(31, 5).Value = '  +1.87%  '

$ws.Cells.Item(32, 5).Value = '  -1.40%  '

$ws.Cells.Item(33, 4).Value = "'7.08"
$ws.Cells.Item(33, 5).Value = '  -2.58%  '

$ws.Cells.Item(34, 5).Value = '  +19.09%  '

$ws.Cells.Item(35, 4).Value = "'28.30"
$ws.Cells.Item(35, 5).Value = '  -2.82%  '

$ws.Cells.Item(36, 4).Value = "'1.00"
$ws.Cells.Item(36, 5).Value = '  +0.31%  '

$ws.Cells.Item(37, 4).Value = "'3.710.47"
$ws.Cells.Item(37, 5).Value = '  -1.38%  '

$ws.Cells.Item(38, 4).Value = "'8.83"
$ws.Cells.Item(38, 5).Value = '  -2.49%  '

$ws.Cells.Item(39, 4).Value = "'0.100"
$ws.Cells.Item(39, 5).Value = '  -1.14%  '

$ws.Cells.Item(40, 4).Value = "'3.22"
$ws.Cells.Item(40, 5).Value = '  -3.95%  '

$ws.Cells.Item(41, 4).Value = "'5.74"
$ws.Cells.Item(41, 5).Value = '  -2.76%  '

$ws.Cells.Item(42, 4).Value = "'0.999"
$ws.Cells.Item(42, 5).Value = '  -0.08%  '

$ws.Cells.Item(43, 4).Value = "'0.955"
$ws.Cells.Item(43, 5).Value = '  -2.59%  '

$ws.Cells.Item(44, 5).Value = '  -0.02%  '

$ws.Cells.Item(45, 2).Value = 'Arweave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(45, 4).Value = "'44.30"
$ws.Cells.Item(45, 5).Value = '  +4.19%  '

$ws.Cells.Item(46, 4).Value = "'155.33"
$ws.Cells.Item(46, 5).Value = '  -1.11%  '

$ws.Cells.Item(47, 2).Value = 'Stacks'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(47, 4).Value = "'1.95"
$ws.Cells.Item(47, 5).Value = '  +3.29%  '

$ws.Cells.Item(48, 4).Value = "'46.85"
$ws.Cells.Item(48, 5).Value = '  +0.26%  '

$ws.Cells.Item(49, 5).Value = '  -3.37%  '

$ws.Cells.Item(50, 5).Value = '  -2.28%  '

$ws.Cells.Item(51, 4).Value = "'8.32"
$ws.Cells.Item(51, 5).Value = '  -1.31%  '
